# Penalty Reward System (unfinished) - data re-shuffle on "Weekly Quantity"
# and two value tweaks on "Monthly Trend".
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

# New Order Week (A) / Requested quantity (B) values for rows 14-53.
$weeklyData = @(
    @(14, "45088.99999999999", 110),
    @(15, "45109.99999999999", 410),
    @(16, "45130.99999999999", 520),
    @(17, "45144.99999999999", 10),
    @(18, "45151.99999999999", 230),
    @(19, "45158.99999999999", 50),
    @(20, "45165.99999999999", 190),
    @(21, "45186.99999999999", 200),
    @(22, "45193.99999999999", 130),
    @(23, "45200.99999999999", 360),
    @(24, "45214.99999999999", 10),
    @(25, "45221.99999999999", 520),
    @(26, "45228.99999999999", 200),
    @(27, "45235.99999999999", 20),
    @(28, "45242.99999999999", 30),
    @(29, "45249.99999999999", 100),
    @(30, "45256.99999999999", 70),
    @(31, "45270.99999999999", 270),
    @(32, "45305.99999999999", 30),
    @(33, "45312.99999999999", 270),
    @(34, "45319.99999999999", 180),
    @(35, "45326.99999999999", 210),
    @(36, "45333.99999999999", 1080),
    @(37, "45347.99999999999", 360),
    @(38, "45361.99999999999", 1940),
    @(39, "45375.99999999999", 310),
    @(40, "45417.99999999999", 80),
    @(41, "45424.99999999999", 60),
    @(42, "45431.99999999999", 20),
    @(43, "45459.99999999999", 20),
    @(44, "45466.99999999999", 40),
    @(45, "45473.99999999999", 40),
    @(46, "45480.99999999999", 60),
    @(47, "45494.99999999999", 980),
    @(48, "45543.99999999999", 80),
    @(49, "45550.99999999999", 100),
    @(50, "45557.99999999999", 100),
    @(51, "45564.99999999999", 120),
    @(52, "45571.99999999999", 200),
    @(53, "45585.99999999999", 300)
)

foreach ($row in $weeklyData) {
    $r = $row[0]
    $wsWeekly.Cells.Item($r, 1).Value = $row[1]
    $wsWeekly.Cells.Item($r, 2).Value = $row[2]
}

# Old rows 54-58 are no longer part of the table - remove them so the
# sheet's used range / dimension shrinks back down to A1:B53.
$wsWeekly.Rows("54:58").Delete()

# --- Sheet 2: "Monthly Trend" ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B5").Value = 520
$wsMonthly.Range("B6").Value = 520
